# "article 91 is live"
# Shift the "ser" (article series) numbers forward by one on row 7:
#   C7: ser 90 -> 91
#   E7: ser 89 -> 90
#   I7: ser 88 -> 89

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 91"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 90"
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 89"
